$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.087.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.258.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  -3.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.827.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '68.082.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.20'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.41%  '
$ws.Range("E16").Value = '  -3.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.268.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '415.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.93%  '
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.506'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.61%  '
$ws.Range("E25").Value = '  -4.16%  '
$ws.Range("E26").Value = '  -1.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.64%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -2.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.63%  '
$ws.Range("E31").Value = '  -6.10%  '
$ws.Range("E32").Value = '  -5.01%  '
$ws.Range("E33").Value = '  -5.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '164.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("E35").Value = '  -5.99%  '
$ws.Range("E36").Value = '  -6.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.789'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.620.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0671'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("E43").Value = '  -5.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '334.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.53%  '
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.976'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0998'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.06%  '
